# Applies the "Saldo" update:
#  1. Remove the LEILA (004208447 / 4000) row.
#  2. Insert a new first data row: 005000645 / ANTONIO / 30109.78
#     right below the header row (i.e. directly above the AHMAD row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: delete the LEILA row (account 004208447) -----------------
# Search column A (by displayed text, so leading zeros compare correctly)
# instead of assuming a fixed row index, in case the sheet layout shifts.
$leilaRow = $null
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp
for ($r = 1; $r -le $lastRow; $r++) {
    if ($ws.Cells.Item($r, 1).Text -eq "004208447") {
        $leilaRow = $r
        break
    }
}
if ($leilaRow -ne $null) {
    $ws.Rows.Item($leilaRow).Delete()
}

# --- Step 2: insert the new ANTONIO row right after the header --------
# The header occupies row 1, so the new row becomes row 2 and everything
# that followed (starting with AHMAD) shifts down by one.
$ws.Rows.Item(2).Insert()

# Write the account number through a scratch cell formatted as Text so the
# leading zeros in "005000645" survive, then copy/paste only the resulting
# value into place and drop the scratch cell - this keeps the destination
# cell on the sheet's normal (unstyled) format, matching the other rows.
$scratch = $ws.Cells.Item($ws.Rows.Count, 50)
$scratch.NumberFormat = "@"
$scratch.Value = "005000645"
$scratch.Copy()
$ws.Cells.Item(2, 1).PasteSpecial(-4163)  # xlPasteValues
$scratch.Clear()
$excel.CutCopyMode = $false

$ws.Cells.Item(2, 2).Value = "ANTONIO"
$ws.Cells.Item(2, 3).Value = 30109.78
